$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was inserted: copy row 35 (whole row, values + formatting)
# down into a newly inserted row 36, pushing rows 36-89 down to 37-90.
$ws.Rows.Item(35).Copy()
$ws.Rows.Item(36).Insert()

# Row 35 now holds the new observation for this week.
$ws.Range("D35").Value = 44571
$ws.Range("J35").Value = 300
